$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize G2 cell style: it previously carried a redundant explicit
# "applyFont" style variant (same font/format as the rest of row 1-3);
# reset it so it shares the common text-format style used elsewhere.
$ws.Range("G2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"

# Re-map the crop/seed layout cells to the values actually seen in-game.
# Every value below already exists in the shared-string table, so this
# just repoints which cell references which existing string.
$ws.Range("G1").Value = "401:10"
$ws.Range("G2").Value = "401:11"
$ws.Range("E3").Value = "401:0"
$ws.Range("F3").Value = "401:5"
$ws.Range("K3").Value = "402:5"
$ws.Range("L3").Value = "402:10"
$ws.Range("A4").Value = "400:0"
$ws.Range("B4").Value = "400:4"
$ws.Range("C4").Value = "400:8"
$ws.Range("D4").Value = "400:12"
$ws.Range("E4").Value = "401:1"
$ws.Range("F4").Value = "401:6"
$ws.Range("G4").Value = "401:13"
$ws.Range("K4").Value = "402:6"
$ws.Range("L4").Value = "402:11"
$ws.Range("M4").Value = "403:0"
$ws.Range("N4").Value = "403:4"
$ws.Range("O4").Value = "403:8"
$ws.Range("P4").Value = "403:12"
$ws.Range("A5").Value = "400:1"
$ws.Range("B5").Value = "400:5"
$ws.Range("C5").Value = "400:9"
$ws.Range("D5").Value = "400:13"
$ws.Range("G5").Value = "401:14"
$ws.Range("K5").Value = "402:7"
$ws.Range("L5").Value = "402:12"
$ws.Range("M5").Value = "403:1"
$ws.Range("N5").Value = "403:5"
$ws.Range("O5").Value = "403:9"
$ws.Range("P5").Value = "403:13"
$ws.Range("A6").Value = "400:2"
$ws.Range("B6").Value = "400:6"
$ws.Range("C6").Value = "400:10"
$ws.Range("D6").Value = "400:14"
$ws.Range("E6").Value = "401:3"
$ws.Range("F6").Value = "401:8"
$ws.Range("K6").Value = "402:8"
$ws.Range("L6").Value = "402:13"
$ws.Range("M6").Value = "403:2"
$ws.Range("N6").Value = "403:6"
$ws.Range("O6").Value = "403:10"
$ws.Range("P6").Value = "403:14"
$ws.Range("A7").Value = "400:3"
$ws.Range("B7").Value = "400:7"
$ws.Range("C7").Value = "400:11"
$ws.Range("D7").Value = "400:15"
$ws.Range("E7").Value = "401:4"
$ws.Range("F7").Value = "401:9"
$ws.Range("K7").Value = "402:9"
$ws.Range("L7").Value = "402:14"
$ws.Range("M7").Value = "403:3"
$ws.Range("N7").Value = "403:7"
$ws.Range("O7").Value = "403:11"
$ws.Range("P7").Value = "403:15"
$ws.Range("A9").Value = "404:0"
$ws.Range("B9").Value = "404:5"
$ws.Range("C9").Value = "404:10"
$ws.Range("D9").Value = "404:15"
$ws.Range("E9").Value = "405:0"
$ws.Range("F9").Value = "405:5"
$ws.Range("K9").Value = "406:10"
$ws.Range("L9").Value = "406:15"
$ws.Range("M9").Value = "407:0"
$ws.Range("N9").Value = "407:5"
$ws.Range("O9").Value = "407:10"
$ws.Range("P9").Value = "407:15"
$ws.Range("A10").Value = "404:1"
$ws.Range("B10").Value = "404:6"
$ws.Range("C10").Value = "404:11"
$ws.Range("D10").Value = "404:16"
$ws.Range("E10").Value = "405:1"
$ws.Range("F10").Value = "405:6"
$ws.Range("K10").Value = "406:11"
$ws.Range("L10").Value = "406:16"
$ws.Range("M10").Value = "407:1"
$ws.Range("N10").Value = "407:6"
$ws.Range("O10").Value = "407:11"
$ws.Range("P10").Value = "407:16"
$ws.Range("E11").Value = "405:2"
$ws.Range("F11").Value = "405:7"
$ws.Range("G11").Value = "405:10"
$ws.Range("H11").Value = "405:15"
$ws.Range("I11").Value = "406:0"
$ws.Range("J11").Value = "406:5"
$ws.Range("K11").Value = "406:12"
$ws.Range("L11").Value = "406:17"
$ws.Range("M11").Value = "407:2"
$ws.Range("N11").Value = "407:7"
$ws.Range("O11").Value = "407:12"
$ws.Range("P11").Value = "407:17"
$ws.Range("A12").Value = "404:3"
$ws.Range("B12").Value = "404:8"
$ws.Range("C12").Value = "404:13"
$ws.Range("D12").Value = "404:18"
$ws.Range("E12").Value = "405:3"
$ws.Range("F12").Value = "405:8"
$ws.Range("G12").Value = "405:11"
$ws.Range("H12").Value = "405:16"
$ws.Range("I12").Value = "406:1"
$ws.Range("J12").Value = "406:6"
$ws.Range("K12").Value = "406:13"
$ws.Range("L12").Value = "406:18"
$ws.Range("M12").Value = "407:3"
$ws.Range("N12").Value = "407:8"
$ws.Range("O12").Value = "407:13"
$ws.Range("P12").Value = "407:18"
$ws.Range("A13").Value = "404:4"
$ws.Range("B13").Value = "404:9"
$ws.Range("C13").Value = "404:14"
$ws.Range("D13").Value = "404:19"
$ws.Range("E13").Value = "405:4"
$ws.Range("F13").Value = "405:9"
$ws.Range("G13").Value = "405:12"
$ws.Range("H13").Value = "405:17"
$ws.Range("I13").Value = "406:2"
$ws.Range("J13").Value = "406:7"
$ws.Range("K13").Value = "406:14"
$ws.Range("L13").Value = "406:19"
$ws.Range("M13").Value = "407:4"
$ws.Range("N13").Value = "407:9"
$ws.Range("O13").Value = "407:14"
$ws.Range("P13").Value = "407:19"
$ws.Range("G14").Value = "405:13"
$ws.Range("H14").Value = "405:18"
$ws.Range("I14").Value = "406:3"
$ws.Range("J14").Value = "406:8"
$ws.Range("G15").Value = "405:14"
$ws.Range("H15").Value = "405:19"
$ws.Range("I15").Value = "406:4"
$ws.Range("J15").Value = "406:9"
